$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 3).Value = 48608
$ws.Cells.Item(2, 4).Value = 82283730
$ws.Cells.Item(3, 3).Value = 116816
$ws.Cells.Item(3, 4).Value = 212122585
$ws.Cells.Item(4, 3).Value = 41089
$ws.Cells.Item(4, 4).Value = 83866911
$ws.Cells.Item(5, 3).Value = 12289
$ws.Cells.Item(5, 4).Value = 27322496
$ws.Cells.Item(6, 3).Value = 3812
$ws.Cells.Item(6, 4).Value = 11037470
$ws.Cells.Item(7, 3).Value = 588
$ws.Cells.Item(7, 4).Value = 2209356
$ws.Cells.Item(10, 3).Value = 6
$ws.Cells.Item(10, 4).Value = 26323
$ws.Cells.Item(12, 3).Value = 50067
$ws.Cells.Item(12, 4).Value = 71748351
$ws.Cells.Item(13, 3).Value = 12347
$ws.Cells.Item(13, 4).Value = 20616898
$ws.Cells.Item(14, 3).Value = 32863
$ws.Cells.Item(14, 4).Value = 57395569
$ws.Cells.Item(15, 3).Value = 10643
$ws.Cells.Item(15, 4).Value = 20231462
$ws.Cells.Item(16, 3).Value = 2906
$ws.Cells.Item(16, 4).Value = 5814288
$ws.Cells.Item(17, 3).Value = 782
$ws.Cells.Item(17, 4).Value = 1996408
$ws.Cells.Item(18, 3).Value = 122
$ws.Cells.Item(18, 4).Value = 407661
$ws.Cells.Item(20, 3).Value = 12396
$ws.Cells.Item(20, 4).Value = 17215770
$ws.Cells.Item(21, 3).Value = 17061
$ws.Cells.Item(21, 4).Value = 28355739
$ws.Cells.Item(22, 3).Value = 40088
$ws.Cells.Item(22, 4).Value = 70126801
$ws.Cells.Item(23, 3).Value = 13433
$ws.Cells.Item(23, 4).Value = 25937599
$ws.Cells.Item(24, 3).Value = 3691
$ws.Cells.Item(24, 4).Value = 7106074
$ws.Cells.Item(25, 3).Value = 977
$ws.Cells.Item(25, 4).Value = 2315862
$ws.Cells.Item(26, 3).Value = 114
$ws.Cells.Item(26, 4).Value = 395554
$ws.Cells.Item(27, 3).Value = 13926
$ws.Cells.Item(27, 4).Value = 19256836
$ws.Cells.Item(28, 3).Value = 9935
$ws.Cells.Item(28, 4).Value = 16958894
$ws.Cells.Item(29, 3).Value = 28511
$ws.Cells.Item(29, 4).Value = 50726391
$ws.Cells.Item(30, 3).Value = 10051
$ws.Cells.Item(30, 4).Value = 19640175
$ws.Cells.Item(31, 3).Value = 2647
$ws.Cells.Item(31, 4).Value = 5258822
$ws.Cells.Item(32, 3).Value = 670
$ws.Cells.Item(32, 4).Value = 1600967
$ws.Cells.Item(34, 3).Value = 10141
$ws.Cells.Item(34, 4).Value = 14126072
$ws.Cells.Item(35, 3).Value = 4318
$ws.Cells.Item(35, 4).Value = 7447443
$ws.Cells.Item(36, 3).Value = 10114
$ws.Cells.Item(36, 4).Value = 18225473
$ws.Cells.Item(37, 3).Value = 4031
$ws.Cells.Item(37, 4).Value = 7942668
$ws.Cells.Item(38, 3).Value = 1062
$ws.Cells.Item(38, 4).Value = 2239144
$ws.Cells.Item(39, 3).Value = 273
$ws.Cells.Item(39, 4).Value = 767012
$ws.Cells.Item(40, 3).Value = 26
$ws.Cells.Item(40, 4).Value = 132075
$ws.Cells.Item(41, 3).Value = 3139
$ws.Cells.Item(41, 4).Value = 4458705
$ws.Cells.Item(42, 3).Value = 21745
$ws.Cells.Item(42, 4).Value = 36216323
$ws.Cells.Item(43, 3).Value = 63973
$ws.Cells.Item(43, 4).Value = 112079257
$ws.Cells.Item(44, 3).Value = 24112
$ws.Cells.Item(44, 4).Value = 45827033
$ws.Cells.Item(45, 3).Value = 7514
$ws.Cells.Item(45, 4).Value = 14542870
$ws.Cells.Item(46, 3).Value = 2182
$ws.Cells.Item(46, 4).Value = 5191053
$ws.Cells.Item(47, 3).Value = 248
$ws.Cells.Item(47, 4).Value = 869659
$ws.Cells.Item(50, 3).Value = 20906
$ws.Cells.Item(50, 4).Value = 30579669
$ws.Cells.Item(51, 3).Value = 2580
$ws.Cells.Item(51, 4).Value = 3992571
$ws.Cells.Item(52, 3).Value = 8622
$ws.Cells.Item(52, 4).Value = 13494779
$ws.Cells.Item(53, 3).Value = 2908
$ws.Cells.Item(53, 4).Value = 4798478
$ws.Cells.Item(54, 3).Value = 943
$ws.Cells.Item(54, 4).Value = 1662564
$ws.Cells.Item(55, 3).Value = 268
$ws.Cells.Item(55, 4).Value = 499005
$ws.Cells.Item(57, 3).Value = 5
$ws.Cells.Item(57, 4).Value = 16000
$ws.Cells.Item(58, 3).Value = 8606
$ws.Cells.Item(58, 4).Value = 12260864
$ws.Cells.Item(59, 3).Value = 1762
$ws.Cells.Item(59, 4).Value = 3772020
$ws.Cells.Item(60, 3).Value = 4168
$ws.Cells.Item(60, 4).Value = 8773890
$ws.Cells.Item(61, 3).Value = 1662
$ws.Cells.Item(61, 4).Value = 3599061
$ws.Cells.Item(62, 3).Value = 557
$ws.Cells.Item(62, 4).Value = 1180012
$ws.Cells.Item(65, 3).Value = 2714
$ws.Cells.Item(65, 4).Value = 5279325
$ws.Cells.Item(66, 3).Value = 20007
$ws.Cells.Item(66, 4).Value = 33603290
$ws.Cells.Item(67, 3).Value = 57582
$ws.Cells.Item(67, 4).Value = 103770224
$ws.Cells.Item(68, 3).Value = 20536
$ws.Cells.Item(68, 4).Value = 40879385
$ws.Cells.Item(69, 3).Value = 6229
$ws.Cells.Item(69, 4).Value = 12761252
$ws.Cells.Item(70, 3).Value = 1747
$ws.Cells.Item(70, 4).Value = 4560537
$ws.Cells.Item(71, 3).Value = 258
$ws.Cells.Item(71, 4).Value = 921278
$ws.Cells.Item(74, 3).Value = 18429
$ws.Cells.Item(74, 4).Value = 25711176
$ws.Cells.Item(75, 3).Value = 74131
$ws.Cells.Item(75, 4).Value = 132067024
$ws.Cells.Item(76, 3).Value = 199894
$ws.Cells.Item(76, 4).Value = 371133086
$ws.Cells.Item(77, 3).Value = 87529
$ws.Cells.Item(77, 4).Value = 181683717
$ws.Cells.Item(78, 3).Value = 30584
$ws.Cells.Item(78, 4).Value = 72704288
$ws.Cells.Item(79, 3).Value = 9916
$ws.Cells.Item(79, 4).Value = 31129393
$ws.Cells.Item(80, 3).Value = 1397
$ws.Cells.Item(80, 4).Value = 6839110
$ws.Cells.Item(81, 3).Value = 80
$ws.Cells.Item(81, 4).Value = 359033
$ws.Cells.Item(86, 3).Value = 70901
$ws.Cells.Item(86, 4).Value = 103979394
$ws.Cells.Item(87, 3).Value = 5427
$ws.Cells.Item(87, 4).Value = 8151075
$ws.Cells.Item(88, 3).Value = 13201
$ws.Cells.Item(88, 4).Value = 20140632
$ws.Cells.Item(89, 3).Value = 4274
$ws.Cells.Item(89, 4).Value = 6631551
$ws.Cells.Item(90, 3).Value = 1512
$ws.Cells.Item(90, 4).Value = 2391008
$ws.Cells.Item(91, 3).Value = 393
$ws.Cells.Item(91, 4).Value = 713012
$ws.Cells.Item(92, 3).Value = 44
$ws.Cells.Item(92, 4).Value = 108725
$ws.Cells.Item(94, 3).Value = 6124
$ws.Cells.Item(94, 4).Value = 8356145
$ws.Cells.Item(95, 3).Value = 2103
$ws.Cells.Item(95, 4).Value = 3463178
$ws.Cells.Item(96, 3).Value = 6594
$ws.Cells.Item(96, 4).Value = 11127309
$ws.Cells.Item(97, 3).Value = 2360
$ws.Cells.Item(97, 4).Value = 4195198
$ws.Cells.Item(98, 3).Value = 877
$ws.Cells.Item(98, 4).Value = 1617006
$ws.Cells.Item(99, 3).Value = 280
$ws.Cells.Item(99, 4).Value = 625631
$ws.Cells.Item(100, 3).Value = 39
$ws.Cells.Item(100, 4).Value = 109164
$ws.Cells.Item(102, 3).Value = 4411
$ws.Cells.Item(102, 4).Value = 6100081
$ws.Cells.Item(103, 3).Value = 985
$ws.Cells.Item(103, 4).Value = 2085998
$ws.Cells.Item(104, 3).Value = 683
$ws.Cells.Item(104, 4).Value = 1565637
$ws.Cells.Item(105, 3).Value = 246
$ws.Cells.Item(105, 4).Value = 548616
$ws.Cells.Item(109, 3).Value = 14209
$ws.Cells.Item(109, 4).Value = 24525223
$ws.Cells.Item(110, 3).Value = 37116
$ws.Cells.Item(110, 4).Value = 66691500
$ws.Cells.Item(111, 3).Value = 12782
$ws.Cells.Item(111, 4).Value = 25088042
$ws.Cells.Item(112, 3).Value = 3693
$ws.Cells.Item(112, 4).Value = 7453680
$ws.Cells.Item(113, 3).Value = 959
$ws.Cells.Item(113, 4).Value = 2393231
$ws.Cells.Item(114, 3).Value = 158
$ws.Cells.Item(114, 4).Value = 558868
$ws.Cells.Item(115, 3).Value = 16
$ws.Cells.Item(115, 4).Value = 46668
$ws.Cells.Item(117, 3).Value = 11794
$ws.Cells.Item(117, 4).Value = 16377315
$ws.Cells.Item(118, 3).Value = 38856
$ws.Cells.Item(118, 4).Value = 64860378
$ws.Cells.Item(119, 3).Value = 82921
$ws.Cells.Item(119, 4).Value = 145553251
$ws.Cells.Item(120, 3).Value = 27309
$ws.Cells.Item(120, 4).Value = 52230960
$ws.Cells.Item(121, 3).Value = 8122
$ws.Cells.Item(121, 4).Value = 16204266
$ws.Cells.Item(122, 3).Value = 2059
$ws.Cells.Item(122, 4).Value = 4795529
$ws.Cells.Item(123, 3).Value = 321
$ws.Cells.Item(123, 4).Value = 1001715
$ws.Cells.Item(124, 3).Value = 19
$ws.Cells.Item(124, 4).Value = 53166
$ws.Cells.Item(127, 3).Value = 30777
$ws.Cells.Item(127, 4).Value = 43185314
$ws.Cells.Item(128, 3).Value = 46836
$ws.Cells.Item(128, 4).Value = 80423409
$ws.Cells.Item(129, 3).Value = 97632
$ws.Cells.Item(129, 4).Value = 173717263
$ws.Cells.Item(130, 3).Value = 30902
$ws.Cells.Item(130, 4).Value = 61527774
$ws.Cells.Item(131, 3).Value = 8759
$ws.Cells.Item(131, 4).Value = 18574522
$ws.Cells.Item(132, 3).Value = 2348
$ws.Cells.Item(132, 4).Value = 6238750
$ws.Cells.Item(133, 3).Value = 271
$ws.Cells.Item(133, 4).Value = 995652
$ws.Cells.Item(136, 3).Value = 38410
$ws.Cells.Item(136, 4).Value = 53503679
$ws.Cells.Item(137, 3).Value = 17078
$ws.Cells.Item(137, 4).Value = 28714764
$ws.Cells.Item(138, 3).Value = 40968
$ws.Cells.Item(138, 4).Value = 72598002
$ws.Cells.Item(139, 3).Value = 14904
$ws.Cells.Item(139, 4).Value = 28494702
$ws.Cells.Item(140, 3).Value = 4089
$ws.Cells.Item(140, 4).Value = 8233228
$ws.Cells.Item(141, 3).Value = 996
$ws.Cells.Item(141, 4).Value = 2482770
$ws.Cells.Item(142, 3).Value = 150
$ws.Cells.Item(142, 4).Value = 478125
$ws.Cells.Item(145, 3).Value = 12978
$ws.Cells.Item(145, 4).Value = 18286273
$ws.Cells.Item(146, 3).Value = 46412
$ws.Cells.Item(146, 4).Value = 81078819
$ws.Cells.Item(147, 3).Value = 107279
$ws.Cells.Item(147, 4).Value = 198545214
$ws.Cells.Item(148, 3).Value = 33334
$ws.Cells.Item(148, 4).Value = 70610137
$ws.Cells.Item(149, 3).Value = 9495
$ws.Cells.Item(149, 4).Value = 22360293
$ws.Cells.Item(150, 3).Value = 2970
$ws.Cells.Item(150, 4).Value = 8992835
$ws.Cells.Item(151, 3).Value = 425
$ws.Cells.Item(151, 4).Value = 1791917
$ws.Cells.Item(152, 3).Value = 32
$ws.Cells.Item(152, 4).Value = 104483
$ws.Cells.Item(153, 3).Value = 35834
$ws.Cells.Item(153, 4).Value = 51405210
